$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 128, shifting existing rows 128:223 down to 129:224
$ws.Rows.Item(128).Insert()

# Populate the new row 128 with the data from the diff
$ws.Cells.Item(128, 1).Value = 7
$ws.Cells.Item(128, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(128, 3).Value = "Ñuble"
$ws.Cells.Item(128, 4).Value = 44673
$ws.Cells.Item(128, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(128, 5).Value = 16
$ws.Cells.Item(128, 6).Value = 100112043
$ws.Cells.Item(128, 7).Value = "Pepino ensalada"
$ws.Cells.Item(128, 8).Value = "Sin especificar"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 30
$ws.Cells.Item(128, 11).Value = 21000
$ws.Cells.Item(128, 12).Value = 22000
$ws.Cells.Item(128, 13).Value = 21500
$ws.Cells.Item(128, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(128, 15).Value = "Región del Maule"
$ws.Cells.Item(128, 16).Value = 269
$ws.Cells.Item(128, 17).Value = 80
$ws.Cells.Item(128, 18).Value = "Hortaliza"
